$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix the duplicate/inconsistent "s1" entry in row 7 (Data_required/scope col B)
# so it matches the canonical "S1" used elsewhere.
$ws.Range("B7").Value = "S1"

# Add the two new columns: step_name (E) and risk_name (F).
# Header row first, then column E data top-to-bottom, then column F data
# top-to-bottom, so new shared strings land in the same order the workbook
# author entered them.
$ws.Range("E1").Value = "step_name"
$ws.Range("F1").Value = "risk_name"

$ws.Range("E2").Value = "wsp1"
$ws.Range("E3").Value = "wsp2"
$ws.Range("E4").Value = "wsp3"
$ws.Range("E5").Value = "wsp4"
$ws.Range("E6").Value = "wsp5"
$ws.Range("E7").Value = "wsp6"

$ws.Range("F2").Value = "r1"
$ws.Range("F3").Value = "r2"
$ws.Range("F4").Value = "r3"
$ws.Range("F5").Value = "r4"
$ws.Range("F6").Value = "r4"
$ws.Range("F7").Value = "r4"

# Match header formatting (bold row style) for the new header cells.
$ws.Range("E1").Font.Bold = $true
$ws.Range("F1").Font.Bold = $true

# Size the new columns to fit their content (closest achievable match to the
# author's manually-sized / autofit widths of 14.57 and 10.14 chars).
$ws.Columns.Item(5).ColumnWidth = 13.59
$ws.Columns.Item(6).ColumnWidth = 9.25

# Update selection to reflect the cell the author landed on after editing.
$ws.Range("B8").Select()
